# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$a1 = $wsHoja1.Range("A1")
$text = $a1.Value()
$text = $text.Replace("1000 Bs = 3.22 = 12352.8 pesos", "1000 Bs = 3.25 = 12460.99 pesos")
$text = $text.Replace("12352.8 pesos = 3.2 = 965.53 Bs", "12460.99 pesos = 3.24 = 971.27 Bs")
$a1.Value = $text

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 307.6
$wsTasas.Range("O10").Value = 3833
$wsTasas.Range("N12").Value = 3848.99
$wsTasas.Range("O12").Value = 300.01
